$d = $word.ActiveDocument

# --- Paragraph 1 (the "**ID__AFFARS_..._ID**" placeholder line):
#     - swap the stale AFFARS topic placeholder for the new MP id
#     - drop the now-orphaned run that held nothing but a trailing space
#     - line its indent/border spacing up with the rest of the body
#       (which already carries this pBdr/225-twip-indent combo)
$p1 = $d.Paragraphs(1)

# 1) Replace the placeholder text.
$find = $p1.Range.Find
$replaced = $find.Execute("**ID__AFFARS_mp_5325_7002_2_topic_3__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP5325_7002_2_2__ID**", 2)

# 2) Remove the trailing run that contains only a single space, sitting
#    right before the paragraph mark.
$pEnd = $p1.Range.End
$trailingSpace = $d.Range($pEnd - 2, $pEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# 3) Match the paragraph formatting used elsewhere in the body: a
#    225-twip (11.25pt) left indent and a paragraph border that only
#    carries spacing (no line), matching the sibling paragraphs below.
$p1.LeftIndent = 11.25

$borders = $p1.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
